# Update MEIC scaling inventory mapping workbook:
#  - last_inv_year sheet: chn (China) last_inv_year bumped from 2012 to 2017
#  - complete_info sheet: new chn/2017/chn row added (keeps the two MEIC
#    inventory versions consistent), shifting later rows down by one
#  - the last_inv_year AutoFilter / _FilterDatabase range is widened to
#    cover all the data rows (A1:B76) instead of just the header row
#  - complete_info becomes the active sheet/selection when the file is
#    saved (mirrors the author's final view state)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("last_inv_year")
$ws2 = $wb.Worksheets.Item("complete_info")

# --- last_inv_year: bump chn's last inventory year 2012 -> 2017 (row 16) ---
$ws1.Range("B16").Value = 2017

# --- complete_info: insert a new chn/2017/chn row right before the cyp
#     rows (old row 31), pushing every following row down by one ---
$ws2.Rows.Item(31).Insert()
$ws2.Range("A31").Value = "chn"
$ws2.Range("B31").Value = 2017
$ws2.Range("C31").Value = "chn"

# --- widen the last_inv_year AutoFilter range to the full data range ---
$ws1.AutoFilterMode = $false
$ws1.Range("A1:B76").AutoFilter()

$filterName = $wb.Names.Item("last_inv_year!_FilterDatabase")
$filterName.RefersTo = "=last_inv_year!`$A`$1:`$B`$76"

# --- final view state: complete_info is the active/selected sheet ---
$ws1.Range("B14").Select()

$ws2.Activate()
$ws2.Range("A30").Select()
